# 966-MS-EI-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-1-LateRepayment-Loanproduct.xlsx
#
# The product name shared string is renamed from
#   "966-MS-EI-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-1-Late Repayment"
# to
#   "966-MS-EI-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-1-LateRepayment"
# (the space before "Repayment" is removed) on both sheets, and the
# active selection on each sheet is moved back to B1 (the product-name
# cell) instead of wherever it had scrolled/landed before.

$wb = $excel.ActiveWorkbook

$newName = "966-MS-EI-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-1-LateRepayment"

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update the product name header cell on the input sheet.
$wsInput.Range("B1").Value = $newName

# Update the matching product name cell on the output sheet.
$wsOutput.Range("B1").Value = $newName

# Restore the selection to B1 on both sheets (and drop the scrolled
# "topLeftCell" view state the input sheet previously had).
$wsInput.Activate()
[void]$wsInput.Range("B1").Select()

$wsOutput.Activate()
[void]$wsOutput.Range("B1").Select()
